$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correção de BUGS: atualizar os UUIDs "esquisitos" das linhas 2 e 3
$ws.Range("A2").Value = "e1156849-0df0-4f1e-90b4-55bf1e1bb753"
$ws.Range("A3").Value = "57a5aaed-c559-421f-998f-75dfeae4188c"

# Remover as linhas 4 a 6 (Rogerinho, Funcionario, funcionario@email.com)
$ws.Rows("4:6").Delete()
